# Issue #30 Should have an return nav bar on add playlist
#
# - Marks issue #30 ("Should have an return nav bar on add playlist") as DONE.
# - Adds two new issues (#31 "styles to css" / message-alert component,
#   #32 "message/alert component" with a description).
# - Extends the used range / AutoFilter / _FilterDatabase name to cover the
#   new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# --- Mark issue #30 (row 30) as DONE --------------------------------------
$ws.Range("C30").Value = "DONE"

# --- Add new issue rows ----------------------------------------------------
# Row 32 -> issue "#31": styles to css
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 2
$ws.Range("E32").Value = "styles to css"

# Row 33 -> issue "#32": message/alert component
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 2
$ws.Range("E33").Value = "message/alert component"
$ws.Range("H33").Value = "based on timed messages and success/failure"
$ws.Rows.Item(33).RowHeight = 29

# --- Refresh the AutoFilter over the new range ------------------------------
# The sheet is already filtered to show only blank Status (column C / field 3)
# rows. Re-apply the same filter over the extended range so the stored
# <autoFilter> ref keeps up with the new data.
$ws.Range("A1:H25").AutoFilter() | Out-Null
$ws.Range("A1:H31").AutoFilter(3, @(""), 7) | Out-Null

# Row 30 now has a Status value (DONE) too, but it should stay visible (its
# hidden state was already stale/out of sync with the filter before this
# edit), so restore it explicitly after the filter recompute hid it.
$ws.Rows.Item(30).Hidden = $false

# Rows 24 & 31 both already carry a Status value and, per the filter, should
# be hidden (matches the recomputed state from above, kept explicit for
# clarity/robustness).
$ws.Rows.Item(24).Hidden = $true
$ws.Rows.Item(31).Hidden = $true

# --- Update the selection to match where the user ended up -----------------
$ws.Range("B33").Select()

# --- Keep the _xlnm._FilterDatabase defined name in sync with the filter ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Issues!_FilterDatabase") {
        $n.RefersTo = "=Issues!`$A`$1:`$H`$31"
    }
}
